# Update "results hierarchichal vs bisg full_sample.xlsx"
# - refresh the "tree: multiplying leafs" confusion-matrix block (rows 15:19)
#   on the "confusion matrix" sheet with the new C/F totals, clearing the
#   now-unused D/E/G cells
# - restore the real D/G counts on the "BISG" block (rows 25:29) that had
#   been zeroed out / mis-summed
# - leave the workbook with the "confusion matrix" tab active, selection on D20

$wb = $excel.ActiveWorkbook

$wsConfusion = $wb.Worksheets.Item("confusion matrix")
$wsMetrics   = $wb.Worksheets.Item("metrics")

# --- "tree: multiplying leafs" block (rows 15-19): new C/F values, clear D/E/G ---
$treeRows = @{
    15 = @{ C = 3853; F = 2001 }
    16 = @{ C = 3131; F = 15531 }
    17 = @{ C = 78;   F = 11688 }
    18 = @{ C = 571;  F = 56745 }
    19 = @{ C = 595;  F = 4835 }
}

foreach ($r in $treeRows.Keys) {
    $vals = $treeRows[$r]
    $wsConfusion.Range("C$r").Value = $vals.C
    $wsConfusion.Range("D$r").ClearContents()
    $wsConfusion.Range("E$r").ClearContents()
    $wsConfusion.Range("F$r").Value = $vals.F
    $wsConfusion.Range("G$r").ClearContents()
}

# --- "BISG" block (rows 25-29): fix D (was 0) and G (was a stale sum) ---
$bisgRows = @{
    25 = @{ D = 902;   G = 1250 }
    26 = @{ D = 16084; G = 1153 }
    27 = @{ D = 349;   G = 2266 }
    28 = @{ D = 4752;  G = 11732 }
    29 = @{ D = 1227;  G = 1353 }
}

foreach ($r in $bisgRows.Keys) {
    $vals = $bisgRows[$r]
    $wsConfusion.Range("D$r").Value = $vals.D
    $wsConfusion.Range("G$r").Value = $vals.G
}

# --- selection / active sheet: "confusion matrix" becomes the active tab,
#     selection moves to D20; "metrics" keeps its D16 selection but is no
#     longer the active tab ---
[void]$wsMetrics.Range("D16").Select()
[void]$wsConfusion.Activate()
[void]$wsConfusion.Range("D20").Select()
